# Fixed issue anthony about Q calculation
# The load power values (p_mw) in columns B:E, rows 2-25 need to be
# rescaled by a factor of 50/49 to reflect the corrected Q calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$factor = 50.0 / 49.0

$firstRow = 2
$lastRow = 25
$firstCol = 2   # column B
$lastCol = 5    # column E

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $old = $cell.Value2
        $cell.Value2 = $old * $factor
    }
}
